# Add the new "2022-Q4" quarter sheet (fund-holding detail) right before the
# existing "2022-Q3" sheet, fill it with its data, and insert the matching
# summary row at the top of the "总计" (totals) sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new sheet, positioned right before "2022-Q3".
# ---------------------------------------------------------------------------
$q3Ref = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($q3Ref)
$newSheet.Name = "2022-Q4"

# Re-fetch worksheet references by name now that the sheet collection has
# shifted - any handle captured before Add() can point at the wrong sheet.
$q4 = $wb.Worksheets.Item("2022-Q4")
$q3 = $wb.Worksheets.Item("2022-Q3")

# Clone the header-row / index-column formatting (bold + border, s="2" style)
# from the neighboring "2022-Q3" sheet so the new sheet matches the look of
# every other quarter sheet, without creating any new style entries.
$q3.Range("A1:H21").Copy()
$q4.Range("A1:H21").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Header row
# ---------------------------------------------------------------------------
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# 3) Data rows (A=index 0-based, B=基金代码, C=基金名称, D=基金规模,
#    E=股票总仓位, F=仓位占比, G=持有市值(亿元), H=仓位排名)
# ---------------------------------------------------------------------------
$rows = @(
  @(0,  "001678", "英大国企改革主题股票",               "4.20", "92.20", "5.48", "0.2302", 9),
  @(1,  "001144", "大成互联网思维混合",                 "2.35", "68.92", "8.05", "0.1892", 3),
  @(2,  "003713", "英大睿盛灵活配置混合A",               "2.39", "93.29", "6.54", "0.1563", 9),
  @(3,  "003714", "英大睿盛灵活配置混合C",               "2.39", "93.29", "6.54", "0.1563", 9),
  @(4,  "014339", "长江智能制造混合A",                   "2.51", "83.09", "4.69", "0.1177", 4),
  @(5,  "180020", "银华成长先锋混合",                   "2.18", "78.61", "3.49", "0.0761", 9),
  @(6,  "004809", "新疆前海联合润丰灵活配置混合A",        "1.23", "90.88", "4.42", "0.0544", 4),
  @(7,  "000458", "英大领先回报混合",                   "1.81", "93.66", "2.78", "0.0503", 3),
  @(8,  "013721", "信澳景气优选混合A",                   "1.08", "92.37", "4.32", "0.0467", 7),
  @(9,  "013346", "富荣信息技术混合C",                   "1.13", "91.06", "3.25", "0.0367", 9),
  @(10, "004890", "中邮健康文娱灵活配置混合",             "0.42", "92.60", "7.16", "0.0301", 2),
  @(11, "013722", "信澳景气优选混合C",                   "0.48", "92.37", "4.32", "0.0207", 7),
  @(12, "014246", "大摩现代服务业混合A",                 "0.17", "86.98", "8.98", "0.0153", 3),
  @(13, "013345", "富荣信息技术混合A",                   "0.45", "91.06", "3.25", "0.0146", 9),
  @(14, "005444", "光大保德信多策略精选18个月定期开放灵活配置混合", "0.54", "29.43", "2.03", "0.0110", 5),
  @(15, "001270", "英大灵活配置混合A",                   "0.29", "92.68", "2.75", "0.0080", 3),
  @(16, "001271", "英大灵活配置混合B",                   "0.28", "92.68", "2.75", "0.0077", 3),
  @(17, "014247", "大摩现代服务业混合C",                 "0.06", "86.98", "8.98", "0.0054", 3),
  @(18, "014340", "长江智能制造混合C",                   "0.09", "83.09", "4.69", "0.0042", 4),
  @(19, "005935", "新疆前海联合润丰灵活配置混合C",        "0.03", "90.88", "4.42", "0.0013", 4)
)

# Columns B, D, E, F, G hold numeric-looking text (fund codes / percentages
# that must keep fixed decimal places and leading zeros), so they are written
# with a leading "'" to force Excel to store them as text instead of numbers
# -- exactly how the source data is typed (t="inlineStr") in every other
# quarter sheet. Column C (fund name) is natural text already; A and H are
# real numbers.
foreach ($r in $rows) {
  $row = [int]$r[0] + 2
  $q4.Cells.Item($row, 1).Value = $r[0]
  $q4.Cells.Item($row, 2).Value = "'" + $r[1]
  $q4.Cells.Item($row, 3).Value = $r[2]
  $q4.Cells.Item($row, 4).Value = "'" + $r[3]
  $q4.Cells.Item($row, 5).Value = "'" + $r[4]
  $q4.Cells.Item($row, 6).Value = "'" + $r[5]
  $q4.Cells.Item($row, 7).Value = "'" + $r[6]
  $q4.Cells.Item($row, 8).Value = $r[7]
}

# ---------------------------------------------------------------------------
# 4) Insert the corresponding "2022-Q4" row at the top of the "总计" sheet,
#    pushing all existing quarters down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Clone formatting (s="2") for the new index-column cell at A10, taken from
# the existing A9 cell, before the values below shift down into it.
$total.Range("A9").Copy()
$total.Range("A10").PasteSpecial(-4122)

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Copy()
$total.Range("B2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 20
$total.Range("D2").Value = 1.23
